$d = $word.ActiveDocument

$replacements = @(
    @("341×8=", "454×5="),
    @("797×9=", "298×5="),
    @("133×4=", "560×6="),
    @("935×4=", "683×2="),
    @("945×5=", "217×4="),
    @("800×3=", "360×9="),
    @("968×6=", "720×5="),
    @("235×5=", "412×4="),
    @("863×2=", "238×6="),
    @("701×9=", "926×6="),
    @("973×5=", "524×9="),
    @("458×2=", "697×5="),
    @("539×9=", "767×6="),
    @("143×5=", "393×8="),
    @("667×2=", "546×7="),
    @("854×8=", "511×3="),
    @("737×9=", "655×4="),
    @("339×2=", "798×4="),
    @("395×8=", "261×6="),
    @("382×8=", "963×9="),
    @("413×8=", "572×2="),
    @("976×5=", "572×8="),
    @("317×2=", "401×6="),
    @("120×7=", "416×4="),
    @("487×3=", "896×6=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
